$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.120168333333333
$ws.Range("H2").Value = 3.360505
$ws.Range("I2").Value = 0.001768092629909379
$ws.Range("J2").Value = 0.001768092629909379
$ws.Range("M2").Value = 721.1356606666667
$ws.Range("N2").Value = 2163.406982
$ws.Range("O2").Value = 0.8508208584329936
$ws.Range("P2").Value = 0.8508208584329936
$ws.Range("Q2").Value = 807.7933311162121
$ws.Range("R2").Value = 7270.139980045909
$ws.Range("S2").Value = 0.001504330089168547
$ws.Range("T2").Value = 0.001504330089168547
$ws.Range("G3").Value = 1.120168333333333
$ws.Range("H3").Value = 3.360505
$ws.Range("I3").Value = 0.001768092629909379
$ws.Range("J3").Value = 0.001768092629909379
$ws.Range("O3").Value = 0.002793596814304166
$ws.Range("P3").Value = 0.002793596814304166
$ws.Range("Q3").Value = 2.652319644088889
$ws.Range("R3").Value = 23.8708767968
$ws.Range("S3").Value = [double]"4.939337938309515E-06"
$ws.Range("T3").Value = [double]"4.939337938309515E-06"
$ws.Range("G4").Value = 1.120168333333333
$ws.Range("H4").Value = 3.360505
$ws.Range("I4").Value = 0.001768092629909379
$ws.Range("J4").Value = 0.001768092629909379
$ws.Range("O4").Value = 0.1463855447527022
$ws.Range("P4").Value = 0.1463855447527022
$ws.Range("Q4").Value = 138.9825668364939
$ws.Range("R4").Value = 1250.843101528445
$ws.Range("S4").Value = 0.0002588232028025222
$ws.Range("T4").Value = 0.0002588232028025222
$ws.Range("I5").Value = 0.9534130698726969
$ws.Range("J5").Value = 0.9534130698726969
$ws.Range("M5").Value = 721.1356606666667
$ws.Range("N5").Value = 2163.406982
$ws.Range("O5").Value = 0.8508208584329936
$ws.Range("P5").Value = 0.8508208584329936
$ws.Range("Q5").Value = 435588.445205879
$ws.Range("R5").Value = 3920296.006852911
$ws.Range("S5").Value = 0.8111837265503237
$ws.Range("T5").Value = 0.8111837265503237
$ws.Range("I6").Value = 0.9534130698726969
$ws.Range("J6").Value = 0.9534130698726969
$ws.Range("O6").Value = 0.002793596814304166
$ws.Range("P6").Value = 0.002793596814304166
$ws.Range("S6").Value = 0.002663451714712321
$ws.Range("T6").Value = 0.002663451714712321
$ws.Range("I7").Value = 0.9534130698726969
$ws.Range("J7").Value = 0.9534130698726969
$ws.Range("O7").Value = 0.1463855447527022
$ws.Range("P7").Value = 0.1463855447527022
$ws.Range("S7").Value = 0.1395658916076608
$ws.Range("T7").Value = 0.1395658916076608
$ws.Range("H8").Value = 85.18441
$ws.Range("I8").Value = 0.04481883749739363
$ws.Range("J8").Value = 0.04481883749739363
$ws.Range("M8").Value = 721.1356606666667
$ws.Range("N8").Value = 2163.406982
$ws.Range("O8").Value = 0.8508208584329936
$ws.Range("P8").Value = 0.8508208584329936
$ws.Range("Q8").Value = 20476.5052612834
$ws.Range("R8").Value = 184288.5473515506
$ws.Range("S8").Value = 0.03813280179350129
$ws.Range("T8").Value = 0.03813280179350129
$ws.Range("H9").Value = 85.18441
$ws.Range("I9").Value = 0.04481883749739363
$ws.Range("J9").Value = 0.04481883749739363
$ws.Range("O9").Value = 0.002793596814304166
$ws.Range("P9").Value = 0.002793596814304166
$ws.Range("Q9").Value = 67.2328367352889
$ws.Range("R9").Value = 605.0955306176
$ws.Range("S9").Value = 0.0001252057616535349
$ws.Range("T9").Value = 0.0001252057616535349
$ws.Range("H10").Value = 85.18441
$ws.Range("I10").Value = 0.04481883749739363
$ws.Range("J10").Value = 0.04481883749739363
$ws.Range("O10").Value = 0.1463855447527022
$ws.Range("P10").Value = 0.1463855447527022
$ws.Range("Q10").Value = 3523.026436875499
$ws.Range("R10").Value = 31707.23793187949
$ws.Range("S10").Value = 0.006560829942238801
$ws.Range("T10").Value = 0.006560829942238802
